# Remove the <w:contextualSpacing w:val="0"/> element from every paragraph's
# pPr. This attribute has no COM-exposed ParagraphFormat property in this
# object model, so we surgically patch each paragraph's OOXML via
# Range.InsertXML (the supported raw-XML escape hatch) rather than trying to
# set a non-existent property.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

$csTag = '<w:contextualSpacing w:val="0"/>'

$patched = 0

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    $xml = $r.WordOpenXML

    if ($xml -notmatch '(?s)<w:body>(.*)</w:body>') {
        continue
    }
    $bodyContent = $matches[1]

    if ($bodyContent -notmatch '(?s)^(<w:p[ >].*?</w:p>)') {
        continue
    }
    $paraXml = $matches[1]

    if ($paraXml -notlike "*$csTag*") {
        # nothing to do for this paragraph
        continue
    }

    $newParaXml = $paraXml.Replace($csTag, '')

    # When the range we fetched is the very last paragraph, Word's exported
    # fragment includes a synthetic trailing empty paragraph that carries the
    # body-level <w:sectPr/> (since the real sectPr lives after our range).
    # If we don't carry that sectPr forward in our replacement, InsertXML
    # leaves a stray extra empty paragraph behind. Detect that tail and graft
    # the real sectPr back on instead of the synthetic paragraph.
    $tail = $bodyContent.Substring($paraXml.Length)
    if ($tail -match '(?s)^<w:p[^>]*/>(<w:sectPr.*)$') {
        $newParaXml = $newParaXml + $matches[1]
    }

    # Collect any relationship ids referenced within the paragraph (e.g.
    # hyperlinks) so the replacement package carries forward the relationships
    # it needs - InsertXML only sees what we hand it.
    $ids = [regex]::Matches($newParaXml, 'r:id="([^"]+)"') |
        ForEach-Object { $_.Groups[1].Value } |
        Select-Object -Unique

    $relsPartXml = ''
    if ($ids -and $ids.Count -gt 0) {
        $relsInner = ''
        if ($xml -match '(?s)pkg:name="/word/_rels/document\.xml\.rels"[^>]*><pkg:xmlData>(.*?)</pkg:xmlData>') {
            $relsXmlAll = $matches[1]
            foreach ($id in $ids) {
                if ($relsXmlAll -match "(<Relationship[^>]*Id=`"$id`"[^>]*/>)") {
                    $relsInner += $matches[1]
                }
            }
        }
        if ($relsInner -ne '') {
            $relsPartXml = '<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' + $relsInner + '</Relationships></pkg:xmlData></pkg:part>'
        }
    }

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
        $newParaXml +
        '</w:body></w:document></pkg:xmlData></pkg:part>' +
        $relsPartXml +
        '</pkg:package>'

    $r.InsertXML($pkg)
    $patched++
}

Write-Host "patched paragraphs: $patched of $count"
